# Logged Week 15 and simulated Week 16
# Update the rolling season-to-date Target Depth totals on both the
# offensive (OFF) and defensive (DEF) sheets for row 2 (the "H" / home total row).

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 369
$wsOff.Range("C2").Value = 259
$wsOff.Range("D2").Value = 95
$wsOff.Range("E2").Value = 44

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 488
$wsDef.Range("C2").Value = 337
$wsDef.Range("D2").Value = 128
$wsDef.Range("E2").Value = 53
